$d = $word.ActiveDocument

function Set-BoldRun([string]$searchText, [int]$skipStart) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $searchText)
        return
    }
    if ($skipStart -gt 0) {
        $rng.MoveStart(1, $skipStart)
    }
    # Re-assert Bold on the run to force the run properties to be
    # rewritten (canonicalizing the <w:b/>/<w:bCs/> element order).
    $rng.Font.Bold = -1
}

Set-BoldRun "(A" 1
Set-BoldRun "A – Anonym" 4
Set-BoldRun "(role B" 6
Set-BoldRun "B – Badatel" 4
Set-BoldRun "(C" 1
Set-BoldRun "C – Archeolog" 4
Set-BoldRun "úrovni D" 7
Set-BoldRun "D – Archivář" 4
Set-BoldRun "oprávněním E" 11
Set-BoldRun "E – Administrátor" 4
